$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in the newly documented component fields (KX224 v2.0 docs update)
$ws.Range("E12").Value = "Pitch 2,54"
$ws.Range("E11").Value = "ADXL355"
$ws.Range("E13").Value = "CMS_0805"
$ws.Range("E14").Value = "CMS_0603 100nF"
$ws.Range("E15").Value = "CMS_0603 10uF"

# Update the selected/active cell to match the latest editor state
$ws.Activate()
$ws.Range("G17").Select()
